# Refresh cryptos list (Price / Volume(1h) columns, and the Stellar/VeChain
# row swap) to match the latest scrape, per the GitHub Actions commit.
# D-column price cells are forced to Text format before assignment (and the
# style is reset back to "Normal" afterwards) so that values such as
# "1.00", "0.0250" or "63.728.03" are stored as literal text instead of
# being auto-converted to numbers/dates by Excel's usual cell-entry rules.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.728.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.735.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.86%  "

$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.219.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.600.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.740.39"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.17%  "

$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("E21").Value = "  -1.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.522"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "

$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

$u1 = [char]0x2083
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0${u1}0905"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.44%  "

$ws.Range("E29").Value = "  +2.30%  "

$ws.Range("E30").Value = "  +5.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("E35").Value = "  +4.09%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.980"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "344.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.04%  "

$ws.Range("E40").Value = "  +1.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0581"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.623"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0250"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
